# Disaggregation of commodity Copper
#
# The "Sector" breakdown on every yearly sheet lists, in row order:
#   row 5 -> Photovoltaic plants
#   row 6 -> Onshore wind plants
#   row 7 -> Offshore wind plants
#
# This edit swaps the "Photovoltaic plants" and "Onshore wind plants"
# rows (both their label in column C and their value in column E) on
# every sheet of the workbook, so the two technologies trade places in
# the listing (their figures travel with their own label).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $label5 = $ws.Range("C5").Value()
    $label6 = $ws.Range("C6").Value()
    $ws.Range("C5").Value = $label6
    $ws.Range("C6").Value = $label5

    $value5 = $ws.Range("E5").Value()
    $value6 = $ws.Range("E6").Value()
    $ws.Range("E5").Value = $value6
    $ws.Range("E6").Value = $value5
}
